# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
# Reverse the "Periodo Mora" (column E, rows 16-39) from descending order
# (2003 .. 1804) to ascending order (1804 .. 2003).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$periodos = @(
    "1804","1805","1806","1807","1808","1809",
    "1810","1811","1812",
    "1901","1902","1903","1904","1905","1906",
    "1907","1908","1909","1910","1911","1912",
    "2001","2002","2003"
)

$startRow = 16
for ($i = 0; $i -lt $periodos.Length; $i++) {
    $row = $startRow + $i
    $ws.Range("E$row").Value = $periodos[$i]
}
